$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text-preserving number format for ID/date-like columns so that
# numeric-looking strings (case numbers, comuna codes, dates-as-text) are
# stored as text, matching the source export, instead of being
# auto-converted to numbers/dates by Excel.
$ws.Range("A20:E38").NumberFormat = "@"

# Row 20
$ws.Range("A20").Value = "804427444"
$ws.Range("B20").Value = "4/1/2025"
$ws.Range("C20").Value = "Cochrane 2864"
$ws.Range("D20").Value = "12"
$ws.Range("E20").Value = "804427444"
$ws.Range("F20").Value = "INCO"
$ws.Range("G20").Value = "Pendiente"
$ws.Range("H20").Value = "Cambiar"
$ws.Range("J20").Value = "Cambio"
$ws.Range("K20").Value = "Sin equipos"
$ws.Range("L20").Value = "Pasante"
$ws.Range("O20").Value = "Paternal"
$ws.Range("P20").Value = "Capital Norte"
$ws.Range("I20").Value = 0
$ws.Range("M20").Value = -58.507569
$ws.Range("N20").Value = -34.579623

# Row 21
$ws.Range("A21").Value = "804568979"
$ws.Range("B21").Value = "4/8/2025"
$ws.Range("C21").Value = "Quesada 2710"
$ws.Range("D21").Value = "13"
$ws.Range("E21").Value = "804568979"
$ws.Range("F21").Value = "INCO"
$ws.Range("G21").Value = "Pendiente"
$ws.Range("H21").Value = "Picada"
$ws.Range("J21").Value = "Cambio"
$ws.Range("K21").Value = "Sin equipos"
$ws.Range("L21").Value = "Pasante"
$ws.Range("O21").Value = "Saavedra"
$ws.Range("P21").Value = "Capital Norte"
$ws.Range("I21").Value = 1
$ws.Range("M21").Value = -58.466348
$ws.Range("N21").Value = -34.556028

# Row 22
$ws.Range("A22").Value = "804736517"
$ws.Range("B22").Value = "4/15/2025"
$ws.Range("C22").Value = "Av. Gral. Mosconi 2490"
$ws.Range("D22").Value = "12"
$ws.Range("E22").Value = "804736517"
$ws.Range("F22").Value = "INCO"
$ws.Range("G22").Value = "Pendiente"
$ws.Range("H22").Value = "Cambiar"
$ws.Range("J22").Value = "Cambio"
$ws.Range("K22").Value = "Sin equipos"
$ws.Range("L22").Value = "Pasante"
$ws.Range("O22").Value = "Paternal"
$ws.Range("P22").Value = "Capital Norte"
$ws.Range("I22").Value = 1
$ws.Range("M22").Value = -58.497446
$ws.Range("N22").Value = -34.583455

# Row 23
$ws.Range("A23").Value = "805507192"
$ws.Range("B23").Value = "4/28/2025"
$ws.Range("C23").Value = "Virrey Arredondo 2821"
$ws.Range("D23").Value = "13"
$ws.Range("E23").Value = "805507192"
$ws.Range("F23").Value = "INCO"
$ws.Range("G23").Value = "Pendiente"
$ws.Range("H23").Value = "Picada"
$ws.Range("J23").Value = "Cambio"
$ws.Range("K23").Value = "Sin equipos"
$ws.Range("L23").Value = "Terminal"
$ws.Range("O23").Value = "Colegiales"
$ws.Range("P23").Value = "Capital Norte"
$ws.Range("I23").Value = 1
$ws.Range("M23").Value = -58.454065
$ws.Range("N23").Value = -34.57105

# Row 24
$ws.Range("A24").Value = "805655355"
$ws.Range("B24").Value = "5/5/2025"
$ws.Range("C24").Value = "Arce 867"
$ws.Range("D24").Value = "14"
$ws.Range("E24").Value = "805655355"
$ws.Range("F24").Value = "INCO"
$ws.Range("G24").Value = "Pendiente"
$ws.Range("H24").Value = "Picada"
$ws.Range("J24").Value = "Cambio"
$ws.Range("K24").Value = "Sin equipos"
$ws.Range("L24").Value = "Pasante"
$ws.Range("O24").Value = "Palermo"
$ws.Range("P24").Value = "Capital Sur"
$ws.Range("I24").Value = 1
$ws.Range("M24").Value = -58.436255
$ws.Range("N24").Value = -34.567733

# Row 25
$ws.Range("A25").Value = "805655369"
$ws.Range("B25").Value = "5/5/2025"
$ws.Range("C25").Value = "Benjamin Matienzo 1524"
$ws.Range("D25").Value = "14"
$ws.Range("E25").Value = "805655369"
$ws.Range("F25").Value = "INCO"
$ws.Range("G25").Value = "Pendiente"
$ws.Range("H25").Value = "Picada"
$ws.Range("J25").Value = "Cambio"
$ws.Range("K25").Value = "Sin equipos"
$ws.Range("L25").Value = "Terminal"
$ws.Range("O25").Value = "Palermo"
$ws.Range("P25").Value = "Capital Sur"
$ws.Range("I25").Value = 1
$ws.Range("M25").Value = -58.43247
$ws.Range("N25").Value = -34.566492

# Row 26
$ws.Range("A26").Value = "6180"
$ws.Range("B26").Value = "5/4/2025"
$ws.Range("C26").Value = "AZARA 15"
$ws.Range("D26").Value = "4"
$ws.Range("E26").Value = "805655333"
$ws.Range("F26").Value = "INCO"
$ws.Range("G26").Value = "Pendiente"
$ws.Range("H26").Value = "Picada"
$ws.Range("J26").Value = "Cambio"
$ws.Range("K26").Value = "Sin equipos"
$ws.Range("L26").Value = "Pasante"
$ws.Range("O26").Value = "San Telmo"
$ws.Range("P26").Value = "Capital Sur"
$ws.Range("I26").Value = 1
$ws.Range("M26").Value = -58.372751
$ws.Range("N26").Value = -34.631917

# Row 27
$ws.Range("A27").Value = "805707245"
$ws.Range("B27").Value = "5/6/2025"
$ws.Range("C27").Value = "Soldado de la Independencia 1298"
$ws.Range("D27").Value = "14"
$ws.Range("E27").Value = "805707245"
$ws.Range("F27").Value = "INCO"
$ws.Range("G27").Value = "Pendiente"
$ws.Range("H27").Value = "Picada - Con fuente teco"
$ws.Range("J27").Value = "Cambio"
$ws.Range("K27").Value = "Fuente Teco"
$ws.Range("L27").Value = "Pasante"
$ws.Range("O27").Value = "Colegiales"
$ws.Range("P27").Value = "Capital Norte"
$ws.Range("I27").Value = 1
$ws.Range("M27").Value = -58.440507
$ws.Range("N27").Value = -34.564016

# Row 28
$ws.Range("A28").Value = "805722772"
$ws.Range("B28").Value = "5/7/2025"
$ws.Range("C28").Value = "Luis Maria Campos 1336"
$ws.Range("D28").Value = "14"
$ws.Range("E28").Value = "805722772"
$ws.Range("F28").Value = "INCO"
$ws.Range("G28").Value = "Pendiente"
$ws.Range("H28").Value = "Picada"
$ws.Range("J28").Value = "Cambio"
$ws.Range("K28").Value = "Sin equipos"
$ws.Range("L28").Value = "Pasante"
$ws.Range("O28").Value = "Colegiales"
$ws.Range("P28").Value = "Capital Norte"
$ws.Range("I28").Value = 1
$ws.Range("M28").Value = -58.44191
$ws.Range("N28").Value = -34.564245

# Row 29
$ws.Range("A29").Value = "-507"
$ws.Range("B29").Value = "7/14/2025"
$ws.Range("C29").Value = "Tamborini 3291"
$ws.Range("D29").Value = "12"
$ws.Range("E29").Value = "808194229"
$ws.Range("F29").Value = "INCO"
$ws.Range("G29").Value = "Pendiente"
$ws.Range("H29").Value = "Picada"
$ws.Range("J29").Value = "Cambio"
$ws.Range("K29").Value = "Sin equipos"
$ws.Range("L29").Value = "Pasante"
$ws.Range("O29").Value = "Saavedra"
$ws.Range("P29").Value = "Capital Norte"
$ws.Range("I29").Value = 1
$ws.Range("M29").Value = -58.473937
$ws.Range("N29").Value = -34.557355

# Row 30
$ws.Range("A30").Value = "-508"
$ws.Range("B30").Value = "7/14/2025"
$ws.Range("C30").Value = "Moldes 2463"
$ws.Range("D30").Value = "12"
$ws.Range("E30").Value = "808194234"
$ws.Range("F30").Value = "INCO"
$ws.Range("G30").Value = "Pendiente"
$ws.Range("H30").Value = "Picada"
$ws.Range("J30").Value = "Cambio"
$ws.Range("K30").Value = "Nodo Teco"
$ws.Range("L30").Value = "Pasante"
$ws.Range("O30").Value = "Saavedra"
$ws.Range("P30").Value = "Capital Norte"
$ws.Range("I30").Value = 1
$ws.Range("M30").Value = -58.462281
$ws.Range("N30").Value = -34.560321

# Row 31
$ws.Range("A31").Value = "-509"
$ws.Range("B31").Value = "7/14/2025"
$ws.Range("C31").Value = "Paso 58"
$ws.Range("D31").Value = "3"
$ws.Range("E31").Value = "808194240"
$ws.Range("F31").Value = "INCO"
$ws.Range("G31").Value = "Pendiente"
$ws.Range("H31").Value = "Picada"
$ws.Range("J31").Value = "Cambio"
$ws.Range("K31").Value = "Sin equipos"
$ws.Range("L31").Value = "Pasante"
$ws.Range("O31").Value = "Almagro"
$ws.Range("P31").Value = "Capital Sur"
$ws.Range("I31").Value = 1
$ws.Range("M31").Value = -58.403422
$ws.Range("N31").Value = -34.609195

# Row 32
$ws.Range("A32").Value = "-510"
$ws.Range("B32").Value = "7/14/2025"
$ws.Range("C32").Value = "Larrea 590"
$ws.Range("D32").Value = "3"
$ws.Range("E32").Value = "808194254"
$ws.Range("F32").Value = "INCO"
$ws.Range("G32").Value = "Pendiente"
$ws.Range("H32").Value = "Picada"
$ws.Range("J32").Value = "Cambio"
$ws.Range("K32").Value = "Fuente Teco"
$ws.Range("L32").Value = "Pasante"
$ws.Range("O32").Value = "Almagro"
$ws.Range("P32").Value = "Capital Sur"
$ws.Range("I32").Value = 1
$ws.Range("M32").Value = -58.402353
$ws.Range("N32").Value = -34.602205

# Row 33
$ws.Range("A33").Value = "-512"
$ws.Range("B33").Value = "7/15/2025"
$ws.Range("C33").Value = "Ciudad de la Paz 3742"
$ws.Range("D33").Value = "12"
$ws.Range("E33").Value = "808240230"
$ws.Range("F33").Value = "INCO"
$ws.Range("G33").Value = "Pendiente"
$ws.Range("H33").Value = "Picada"
$ws.Range("J33").Value = "Cambio"
$ws.Range("K33").Value = "Sin equipos"
$ws.Range("L33").Value = "Pasante"
$ws.Range("O33").Value = "Saavedra"
$ws.Range("P33").Value = "Capital Norte"
$ws.Range("I33").Value = 1
$ws.Range("M33").Value = -58.470347
$ws.Range("N33").Value = -34.547965

# Row 34
$ws.Range("A34").Value = "-514"
$ws.Range("B34").Value = "7/15/2025"
$ws.Range("C34").Value = "Bilbao 2452"
$ws.Range("D34").Value = "7"
$ws.Range("E34").Value = "808243829"
$ws.Range("F34").Value = "INCO"
$ws.Range("G34").Value = "Pendiente"
$ws.Range("H34").Value = "Picada"
$ws.Range("J34").Value = "Cambio"
$ws.Range("K34").Value = "Sin equipos"
$ws.Range("L34").Value = "Pasante"
$ws.Range("O34").Value = "Boedo"
$ws.Range("P34").Value = "Capital Sur"
$ws.Range("I34").Value = 1
$ws.Range("M34").Value = -58.460594
$ws.Range("N34").Value = -34.635581

# Row 35
$ws.Range("A35").Value = "-516"
$ws.Range("B35").Value = "7/16/2025"
$ws.Range("C35").Value = "Olazabal 4417"
$ws.Range("D35").Value = "12"
$ws.Range("E35").Value = "808373646"
$ws.Range("F35").Value = "INCO"
$ws.Range("G35").Value = "Pendiente"
$ws.Range("H35").Value = "Picada"
$ws.Range("J35").Value = "Cambio"
$ws.Range("K35").Value = "Sin equipos"
$ws.Range("L35").Value = "Pasante"
$ws.Range("O35").Value = "Colegiales"
$ws.Range("P35").Value = "Capital Norte"
$ws.Range("I35").Value = 1
$ws.Range("M35").Value = -58.478941
$ws.Range("N35").Value = -34.57242

# Row 36
$ws.Range("A36").Value = "-522"
$ws.Range("B36").Value = "7/21/2025"
$ws.Range("C36").Value = "Uruguay 1090"
$ws.Range("D36").Value = "2"
$ws.Range("E36").Value = "808430941"
$ws.Range("F36").Value = "INCO"
$ws.Range("G36").Value = "Pendiente"
$ws.Range("H36").Value = "Reclaman columna corroida y rienda fuera de norma pero no se ve en la foto."
$ws.Range("J36").Value = "Cambio"
$ws.Range("K36").Value = "Sin equipos"
$ws.Range("L36").Value = "Terminal"
$ws.Range("O36").Value = "Recoleta"
$ws.Range("P36").Value = "Capital Sur"
$ws.Range("I36").Value = 1
$ws.Range("M36").Value = -58.387175
$ws.Range("N36").Value = -34.596

# Row 37
$ws.Range("A37").Value = "-523"
$ws.Range("B37").Value = "7/20/2025"
$ws.Range("C37").Value = "Luis Maria Campos 585"
$ws.Range("D37").Value = "14"
$ws.Range("E37").Value = "808460898"
$ws.Range("F37").Value = "INCO"
$ws.Range("G37").Value = "Pendiente"
$ws.Range("H37").Value = "Picada"
$ws.Range("J37").Value = "Cambio"
$ws.Range("K37").Value = "Sin equipos"
$ws.Range("L37").Value = "Pasante"
$ws.Range("O37").Value = "Palermo"
$ws.Range("P37").Value = "Capital Sur"
$ws.Range("I37").Value = 1
$ws.Range("M37").Value = -58.434668
$ws.Range("N37").Value = -34.571258

# Row 38
$ws.Range("A38").Value = "-524"
$ws.Range("B38").Value = "7/21/2025"
$ws.Range("C38").Value = "Luis Maria Campos 509"
$ws.Range("D38").Value = "14"
$ws.Range("E38").Value = "808460897"
$ws.Range("F38").Value = "INCO"
$ws.Range("G38").Value = "Pendiente"
$ws.Range("H38").Value = "Picada"
$ws.Range("J38").Value = "Cambio"
$ws.Range("K38").Value = "Sin equipos"
$ws.Range("L38").Value = "Pasante"
$ws.Range("O38").Value = "Palermo"
$ws.Range("P38").Value = "Capital Sur"
$ws.Range("I38").Value = 1
$ws.Range("M38").Value = -58.434194
$ws.Range("N38").Value = -34.571754

